$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" = strikeouts) values per regenerated save_data
$kValues = @{
    2 = 2
    3 = 3
    4 = 0
    5 = 1
    6 = 0
    7 = 0
    8 = 0
    9 = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 3
    15 = 1
    16 = 0
    17 = 1
    18 = 3
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 2
    34 = 3
    35 = 2
    36 = 2
    37 = 2
    38 = 0
    39 = 1
    40 = 3
    41 = 0
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 0
    48 = 1
    49 = 0
    50 = 4
    51 = 2
    52 = 0
    53 = 1
    54 = 0
    55 = 1
    56 = 1
    57 = 0
    58 = 1
    59 = 1
    60 = 1
    62 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
